$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the "location" key to "locations" and the "date" key to "dates"
$ws.Range("A2").Value = "locations"
$ws.Range("A3").Value = "dates"

# Update the active selection to A3 (as captured when the sheet was saved)
$ws.Range("A3").Select()
